# Updates the cryptocurrency price table (Sheet1!B2:E51) to the latest
# scrape: the coin list shifted down one row (OKB newly inserted at row 8,
# pushing Cardano/Dogecoin/TRON/... down by one), and every Price (column D)
# and Volume(1h) (column E) figure was refreshed. Column A (rank index) is
# untouched.
#
# Cells are written through a small helper that forces the "@" (text)
# number format before assigning the value and resets the style back to
# "Normal" afterwards. Without this, Excel's COM layer auto-coerces
# plain-looking numeric strings (e.g. "39.73") into real numbers, which
# would change the cell's stored type from text to numeric and diverge
# from the source data (every Price/Volume cell in this sheet is stored as
# text, including values that look numeric, e.g. "25.042.20").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($Worksheet, $Address, $Text)
    $cell = $Worksheet.Range($Address)
    $cell.NumberFormat = "@"
    $cell.Value = $Text
    $cell.Style = "Normal"
}

$updates = [ordered]@{
    # Row 2
    'D2' = '25.029.82'
    'E2' = '  -3.16%  '
    # Row 3
    'D3' = '1.650.63'
    'E3' = '  -4.89%  '
    # Row 4
    'E4' = '  +0.08%  '
    # Row 5
    'D5' = '236.60'
    'E5' = '  -1.83%  '
    # Row 7
    'D7' = '0.4785'
    'E7' = '  -7.83%  '
    # Row 8
    'B8' = 'OKB'
    'C8' = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
    'D8' = '39.73'
    'E8' = '  -0.27%  '
    # Row 9
    'B9' = 'Cardano'
    'C9' = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
    'D9' = '0.2613'
    'E9' = '  -4.34%  '
    # Row 10
    'B10' = 'Dogecoin'
    'C10' = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
    'D10' = '0.05992'
    'E10' = '  -2.48%  '
    # Row 11
    'B11' = 'TRON'
    'C11' = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
    'D11' = '0.07085'
    'E11' = '  -1.13%  '
    # Row 12
    'B12' = 'WrappedEther'
    'C12' = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
    'D12' = '1.656.09'
    'E12' = '  -4.57%  '
    # Row 13
    'B13' = 'Solana'
    'C13' = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
    'D13' = '14.41'
    'E13' = '  -3.44%  '
    # Row 14
    'B14' = 'Polygon'
    'C14' = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
    'D14' = '0.6165'
    'E14' = '  -3.52%  '
    # Row 15
    'B15' = 'Polkadot'
    'C15' = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
    'D15' = '4.569'
    'E15' = '  -0.67%  '
    # Row 16
    'B16' = 'Litecoin'
    'C16' = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
    'D16' = '72.99'
    'E16' = '  -5.20%  '
    # Row 17
    'B17' = 'Dai'
    'C17' = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
    'D17' = '1.001'
    'E17' = '  +0.04%  '
    # Row 18
    'B18' = 'BinanceUSD'
    'C18' = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
    'D18' = '1.001'
    'E18' = '  +0.06%  '
    # Row 19
    'B19' = 'WrappedBTC'
    'C19' = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
    'D19' = '25.026.27'
    'E19' = '  -3.29%  '
    # Row 20
    'B20' = 'Avalanche'
    'C20' = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
    'D20' = '11.33'
    'E20' = '  -3.25%  '
    # Row 21
    'B21' = 'ShibaInu'
    'C21' = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
    'D21' = '0.000006565'
    'E21' = '  -2.89%  '
    # Row 22
    'B22' = 'Uniswap'
    'C22' = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
    'D22' = '4.415'
    'E22' = '  +3.49%  '
    # Row 23
    'B23' = 'WrappedliquidstakedEther2.0'
    'C23' = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
    'D23' = '1.869.31'
    'E23' = '  -4.76%  '
    # Row 24
    'B24' = 'Cosmos'
    'C24' = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
    'D24' = '8.425'
    'E24' = '  -1.99%  '
    # Row 25
    'B25' = 'Chainlink'
    'C25' = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
    'D25' = '5.228'
    'E25' = '  -0.51%  '
    # Row 26
    'B26' = 'Monero'
    'C26' = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
    'D26' = '133.24'
    'E26' = '  -3.17%  '
    # Row 27
    'B27' = 'EthereumClassic'
    'C27' = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
    'D27' = '14.70'
    'E27' = '  -3.10%  '
    # Row 28
    'B28' = 'Toncoin'
    'C28' = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
    'D28' = '1.391'
    'E28' = '  -8.10%  '
    # Row 29
    'B29' = 'LidoDAOToken'
    'C29' = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
    'D29' = '1.683'
    'E29' = '  -4.49%  '
    # Row 30
    'B30' = 'BitcoinCash'
    'C30' = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
    'D30' = '101.40'
    'E30' = '  -3.32%  '
    # Row 31
    'B31' = 'InternetComputer(DFINITY)'
    'C31' = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
    'D31' = '3.763'
    'E31' = '  -4.10%  '
    # Row 32
    'B32' = 'Stellar'
    'C32' = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
    'D32' = '0.07908'
    'E32' = '  -4.06%  '
    # Row 33
    'B33' = 'Filecoin'
    'C33' = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
    'D33' = '3.516'
    'E33' = '  -3.49%  '
    # Row 34
    'B34' = 'Hedera'
    'C34' = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
    'D34' = '0.04551'
    'E34' = '  -1.59%  '
    # Row 35
    'B35' = 'HuobiToken'
    'C35' = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
    'D35' = '2.609'
    'E35' = '  -1.31%  '
    # Row 36
    'B36' = 'ARBITRUM'
    'C36' = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
    'D36' = '0.9358'
    'E36' = '  -5.01%  '
    # Row 37
    'B37' = 'ImmutableX'
    'C37' = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
    'D37' = '0.5771'
    'E37' = '  -6.44%  '
    # Row 38
    'B38' = 'MXToken'
    'C38' = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
    'D38' = '2.618'
    'E38' = '  -2.42%  '
    # Row 39
    'D39' = '0.01533'
    'E39' = '  -3.84%  '
    # Row 40
    'B40' = 'TrustWalletToken'
    'C40' = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
    'D40' = '0.8406'
    'E40' = '  +12.71%  '
    # Row 41
    'B41' = 'PaxDollar'
    'C41' = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
    'D41' = '1.001'
    'E41' = '  +0.11%  '
    # Row 42
    'B42' = 'RenderToken'
    'C42' = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
    'D42' = '1.824'
    'E42' = '  -4.77%  '
    # Row 43
    'B43' = 'Quant'
    'C43' = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
    'D43' = '98.66'
    'E43' = '  -1.39%  '
    # Row 44
    'B44' = 'TheSandbox'
    'C44' = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
    'D44' = '0.3687'
    'E44' = '  -3.79%  '
    # Row 45
    'B45' = 'FraxShare'
    'C45' = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
    'D45' = '4.815'
    'E45' = '  -3.56%  '
    # Row 46
    'B46' = 'Algorand'
    'C46' = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
    'D46' = '0.1108'
    'E46' = '  -1.18%  '
    # Row 47
    'B47' = 'Aptos'
    'C47' = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
    'D47' = '6.011'
    'E47' = '  -3.41%  '
    # Row 48
    'B48' = 'Cronos'
    'C48' = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
    'D48' = '0.05148'
    'E48' = '  -1.66%  '
    # Row 49
    'B49' = 'Aave'
    'C49' = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
    'D49' = '52.03'
    'E49' = '  -4.95%  '
    # Row 50
    'B50' = 'Elrond'
    'C50' = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
    'D50' = '29.28'
    'E50' = '  -3.99%  '
    # Row 51
    'B51' = 'TrueUSD'
    'C51' = 'https://coinranking.com/coin/1ZZI6g5k5royD+trueusd-tusd'
    'D51' = '1.001'
}

foreach ($address in $updates.Keys) {
    Set-TextCell $ws $address $updates[$address]
}

Write-Host "Updated $($updates.Count) cells in the cryptos price table"
